# New weekly record for "Haba" at Feria Lagunitas de Puerto Montt.
# Insert a new row above row 17, shifting the existing rows 17-48 down to 18-49,
# then populate the new row 17 with this week's market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("17:17").Insert()

$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = 44495
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112026
$ws.Range("G17").Value = "Haba"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 10000
$ws.Range("N17").Value = '$/saco 25 kilos'
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 400
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
